$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("List")

# Fill in the newly added API rows (18-21) with key + description.
# Keys (column B) are entered first, then descriptions (column C), so that
# the shared-string table is populated in the same order as the source
# workbook (all keys, then all descriptions).
$ws.Range("B18").Value = "transaction.delete.master.setCountryAdministrativeAreaLevel1"
$ws.Range("B19").Value = "transaction.delete.master.setCountryAdministrativeAreaLevel2"
$ws.Range("B20").Value = "transaction.delete.master.setCountryAdministrativeAreaLevel3"
$ws.Range("B21").Value = "transaction.delete.master.setCountryAdministrativeAreaLevel4"

$ws.Range("C18").Value = "Menghapusi Data Propinsi (Daerah Tingkat 1)"
$ws.Range("C19").Value = "Menghapusi Data Kabupaten / Kota (Daerah Tingkat 2)"
$ws.Range("C20").Value = "Menghapusi Data Kecamatan (Daerah Tingkat 3)"
$ws.Range("C21").Value = "Menghapusi Data Kelurahan / Desa (Daerah Tingkat 4)"

# Update the frozen-pane view so the newly added rows are visible
$ws.Activate()
$ws.Range("C29").Select()
$excel.ActiveWindow.ScrollRow = 12
